# "delivery boy gets notification and information modify"
#
# At the end of the document there are 7 trailing empty paragraphs
# (right after the closing </w:tbl>, before the final sectPr). The edit:
#   1. Gives the 2nd of those 7 (counting from the table) the text
#      "EXTRA:"
#   2. Inserts a brand-new, still-empty paragraph right after it
#   3. Inserts another new, still-empty paragraph right after that one
#   4. Fills that last new paragraph with
#      "DELIVERY BOY GETS NOTIFICATION EMAIL    -DONE"
#   5. Leaves the remaining trailing empty paragraphs untouched, in
#      their original order.
#
# Notes on this runtime's quirks (discovered empirically):
#   - Touching $d.Tables at all corrupts Range.Start/End for every
#     paragraph fetched afterwards in the same script, so table access
#     is avoided entirely; navigation instead walks back from
#     Paragraphs.Last via .Previous().
#   - Range.LanguageID only reliably produces <w:rPr><w:lang .../></w:rPr>
#     on the inserted run when set on the *whole paragraph range*
#     (paragraph.Range, which includes its own end-of-paragraph mark),
#     re-fetched fresh after the text insertion - not on a manually
#     sliced sub-range.
#   - To split a paragraph into "itself, then a new paragraph after it"
#     without leaving stray empty runs, insert a bare CR
#     ([char]13) at the position of the paragraph's own mark
#     (Range.End - 1), rather than using InsertParagraphAfter (which
#     leaves a trailing empty <w:r><w:rPr>.../></w:r>) or Range.Start
#     (which, once the paragraph has text, sits before that text and
#     so splits off the text into the *new* paragraph instead).

$d = $word.ActiveDocument

function GetParaFromEnd($stepsBack) {
    $p = $d.Paragraphs.Last
    for ($k = 0; $k -lt $stepsBack; $k++) {
        $p = $p.Previous()
    }
    return $p
}

# --- Step 1: add "EXTRA:" run to the 2nd of the 7 trailing paragraphs -
# (Last, -1, -2, -3, -4, -5 => 5 steps back from Last.)
$pExtra = GetParaFromEnd 5
$pExtra.Range.InsertAfter("EXTRA:")
$pExtra2 = GetParaFromEnd 5
$pExtra2.Range.LanguageID = "en-GB"

# --- Step 2: insert a brand new, empty paragraph right after it -------
$pExtra3 = GetParaFromEnd 5
$markPos = $pExtra3.Range.End - 1
$d.Range($markPos, $markPos).Text = [char]13

# --- Step 3: insert another new, empty paragraph right after that one -
$pNewEmpty = GetParaFromEnd 5
$markPos2 = $pNewEmpty.Range.End - 1
$d.Range($markPos2, $markPos2).Text = [char]13

# --- Step 4: fill the newest paragraph with the delivery-boy text -----
$pDelivery = GetParaFromEnd 5
$pDelivery.Range.InsertAfter("DELIVERY BOY GETS NOTIFICATION EMAIL    -DONE")
$pDelivery2 = GetParaFromEnd 5
$pDelivery2.Range.LanguageID = "en-GB"
